$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking strings
# (e.g. "1.05") are preserved as text, matching the source data which
# stores all Price/Volume values as inline strings, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '56.984.46'
$ws.Range('E2').Value = '  +10.99%  '
$ws.Range('D3').Value = '3.265.60'
$ws.Range('E3').Value = '  +6.32%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '399.43'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('D6').Value = '110.29'
$ws.Range('E6').Value = '  +8.38%  '
$ws.Range('E7').Value = '  +4.65%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '0.623'
$ws.Range('E9').Value = '  +6.55%  '
$ws.Range('D10').Value = '39.45'
$ws.Range('E10').Value = '  +6.71%  '
$ws.Range('D11').Value = '0.0953'
$ws.Range('E11').Value = '  +12.21%  '
$ws.Range('E12').Value = '  +2.26%  '
$ws.Range('D13').Value = '3.775.90'
$ws.Range('E13').Value = '  +6.44%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '8.12'
$ws.Range('E14').Value = '  +5.63%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '19.13'
$ws.Range('E15').Value = '  +4.30%  '
$ws.Range('D16').Value = '3.259.58'
$ws.Range('E16').Value = '  +6.23%  '
$ws.Range('D17').Value = '1.05'
$ws.Range('D18').Value = '10.89'
$ws.Range('E18').Value = '  +3.08%  '
$ws.Range('D19').Value = '56.807.23'
$ws.Range('E19').Value = '  +10.68%  '
$ws.Range('D20').Value = '3.30'
$ws.Range('E20').Value = '  +4.31%  '
$ws.Range('D21').Value = '0.0000107'
$ws.Range('E21').Value = '  +12.08%  '
$ws.Range('D22').Value = '12.94'
$ws.Range('E22').Value = '  +5.05%  '
$ws.Range('D23').Value = '305.05'
$ws.Range('E23').Value = '  +15.19%  '
$ws.Range('D24').Value = '75.08'
$ws.Range('E24').Value = '  +7.46%  '
$ws.Range('D25').Value = '3.19'
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('D26').Value = '8.02'
$ws.Range('E26').Value = '  +1.17%  '
$ws.Range('D27').Value = '28.29'
$ws.Range('E27').Value = '  +5.06%  '
$ws.Range('D28').Value = '4.38'
$ws.Range('E28').Value = '  +5.25%  '
$ws.Range('D29').Value = '7.29'
$ws.Range('E29').Value = '  +2.09%  '
$ws.Range('D30').Value = '0.169'
$ws.Range('E30').Value = '  +4.56%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  +4.79%  '
$ws.Range('D33').Value = '11.03'
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('D34').Value = '37.81'
$ws.Range('E34').Value = '  +5.14%  '
$ws.Range('D35').Value = '0.0485'
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('E36').Value = '  +3.30%  '
$ws.Range('D37').Value = '51.63'
$ws.Range('E37').Value = '  +3.37%  '
$ws.Range('D38').Value = '3.18'
$ws.Range('E38').Value = '  +26.42%  '
$ws.Range('D39').Value = '3.57'
$ws.Range('E39').Value = '  +7.30%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').Value = '17.54'
$ws.Range('E41').Value = '  +5.13%  '
$ws.Range('E42').Value = '  +5.33%  '
$ws.Range('D43').Value = '134.04'
$ws.Range('E43').Value = '  +4.99%  '
$ws.Range('D44').Value = '4.00'
$ws.Range('E44').Value = '  +2.16%  '
$ws.Range('E45').Value = '  +4.29%  '
$ws.Range('D46').Value = '0.282'
$ws.Range('E46').Value = '  -3.82%  '
$ws.Range('D47').Value = '22.14'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('D48').Value = '2.147.36'
$ws.Range('E48').Value = '  +3.73%  '
$ws.Range('E49').Value = '  +2.51%  '
$ws.Range('D50').Value = '2.39'
$ws.Range('E50').Value = '  -3.25%  '
$ws.Range('D51').Value = '2.01'
$ws.Range('E51').Value = '  +39.89%  '

# Restore the default style on column D so no stray formatting remains
$ws.Range("D2:D51").Style = "Normal"
